$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.844.20'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '3.491.32'
$ws.Range("E3").Value = '  +2.22%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.996'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.35%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '558.14'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.14%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '181.48'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.49%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.638'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +4.42%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.19%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.639'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.74%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.154'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +3.97%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '54.67'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.40%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0000274'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.57%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '9.33'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").Value = '3.989.33'
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '18.78'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +3.15%  '
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.439.07'
$ws.Range("E17").Value = '  +1.11%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '12.11'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +3.87%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '65.569.49'
$ws.Range("E19").Value = '  -1.33%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.996'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.06%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '420.14'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +3.43%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.07'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +5.31%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '86.74'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +3.24%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '4.14'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.81%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '12.79'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +9.54%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '10.94'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -8.64%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.89'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.41%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '6.03'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -3.03%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.14'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +6.42%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '30.56'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.84%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.65'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.57%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '612.52'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -8.35%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '11.84'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.35%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.111'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.82%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '59.57'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '37.91'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.39%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.22%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.146'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +10.12%  '
$ws.Range("D39").Value = '0.0{0}0799' -f [char]0x2083
$ws.Range("E39").Value = '  -2.11%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '3.319.11'
$ws.Range("E40").Value = '  +9.48%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.385'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -4.26%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '3.33'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.994'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.68%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.85'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.29%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.57'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -8.26%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0417'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '3.27'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.06%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.72'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.14%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.133'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.85%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '138.32'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.89%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '8.46'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.68%  '
